$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.716.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.894.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.01%  '
$ws.Range("E4").Value = '  -1.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4878'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3793'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07323'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9129'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.58'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07644'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.901.13'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.479'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.625'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.35'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008768'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.15%  '
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.728.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.52%  '
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.137.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.87%  '
$ws.Range("E24").Value = '  -1.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.889'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.06%  '
$ws.Range("E27").Value = '  -1.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.157'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.885'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08905'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.201'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.226'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7684'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.635'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.568'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02039'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("E38").Value = '  -2.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05284'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5474'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.983'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.886'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.513'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '112.66'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1519'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.66'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4789'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.639'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06051'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.07%  '
